$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2684.375
$ws.Range("I2").Value = 651.2857
$ws.Range("J2").Value = 4265.6665
$ws.Range("K2").Value = 651.2857
$ws.Range("L2").Value = 4265.6665
$ws.Range("M2").Value = -538.2857
$ws.Range("N2").Value = -4491.6665
$ws.Range("H9").Value = 269.8
$ws.Range("I9").Value = 50
$ws.Range("J9").Value = 416.33334
$ws.Range("K9").Value = 50
$ws.Range("L9").Value = 416.33334
$ws.Range("M9").Value = 119
$ws.Range("N9").Value = -754.33334
$ws.Range("H19").Value = 1321.4286
$ws.Range("I19").Value = 1170
$ws.Range("K19").Value = 1170
$ws.Range("M19").Value = -995
$ws.Range("H64").Value = 5272.5454
$ws.Range("I64").Value = 4999.8
$ws.Range("J64").Value = 8000
$ws.Range("K64").Value = 4999.8
$ws.Range("L64").Value = 8000
$ws.Range("M64").Value = -4751.8
$ws.Range("N64").Value = -8496
$ws.Range("H67").Value = 5272.5454
$ws.Range("I67").Value = 4999.8
$ws.Range("J67").Value = 8000
$ws.Range("K67").Value = 4999.8
$ws.Range("L67").Value = 8000
$ws.Range("M67").Value = -4141.8
$ws.Range("N67").Value = -9716
$ws.Range("H70").Value = 2625.45
$ws.Range("I70").Value = 1322.7
$ws.Range("J70").Value = 3928.2
$ws.Range("K70").Value = 3968.1
$ws.Range("L70").Value = 11784.6
$ws.Range("M70").Value = -3698.1
$ws.Range("N70").Value = -12324.6
$ws.Range("H73").Value = 2625.45
$ws.Range("I73").Value = 1322.7
$ws.Range("J73").Value = 3928.2
$ws.Range("K73").Value = 3968.1
$ws.Range("L73").Value = 11784.6
$ws.Range("M73").Value = -3032.1
$ws.Range("N73").Value = -13656.6
$ws.Range("H112").Value = 1633.1578
$ws.Range("I112").Value = 1200
$ws.Range("J112").Value = 1714.375
$ws.Range("K112").Value = 3600
$ws.Range("L112").Value = 5143.125
$ws.Range("M112").Value = -2492
$ws.Range("N112").Value = -7359.125
$ws.Range("H125").Value = 4009.875
$ws.Range("I125").Value = 2813.25
$ws.Range("K125").Value = 25319.25
$ws.Range("M125").Value = -22859.25
$ws.Range("H137").Value = 2513.8333
$ws.Range("I137").Value = 1524.5
$ws.Range("J137").Value = 4492.5
$ws.Range("K137").Value = 4573.5
$ws.Range("L137").Value = 13477.5
$ws.Range("M137").Value = -2023.5
$ws.Range("N137").Value = -18577.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H42").Value = 20000
$ws.Range("J42").Value = 20000
$ws.Range("L42").Value = 20000
$ws.Range("N42").Value = -20972
$ws.Range("H61").Value = 3749.5
$ws.Range("I61").Value = 4000
$ws.Range("K61").Value = 4000
$ws.Range("M61").Value = -3788
$ws.Range("H63").Value = 3283.3333
$ws.Range("I63").Value = 2566.6667
$ws.Range("K63").Value = 2566.6667
$ws.Range("M63").Value = -1880.6667
$ws.Range("H66").Value = 3283.3333
$ws.Range("I66").Value = 2566.6667
$ws.Range("K66").Value = 12833.3335
$ws.Range("M66").Value = -9401.333500000001
$ws.Range("H97").Value = 788.2143
$ws.Range("I97").Value = 748.8461
$ws.Range("J97").Value = 1300
$ws.Range("K97").Value = 748.8461
$ws.Range("L97").Value = 1300
$ws.Range("M97").Value = -252.8461
$ws.Range("N97").Value = -2292
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("N98").Value = ""
$ws.Range("H114").Value = 60000
$ws.Range("J114").Value = 60000
$ws.Range("L114").Value = 60000
$ws.Range("N114").Value = -68678
$ws.Range("H132").Value = 1956.6666
$ws.Range("J132").Value = 2612.375
$ws.Range("L132").Value = 7837.125
$ws.Range("N132").Value = -12897.125
$ws.Range("H136").Value = 3749.5
$ws.Range("I136").Value = 4000
$ws.Range("K136").Value = 12000
$ws.Range("M136").Value = -9450

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 42500
$ws.Range("I102").Value = 42500
$ws.Range("K102").Value = 42500
$ws.Range("M102").Value = -39255
$ws.Range("H105").Value = 2300
$ws.Range("I105").Value = 2200
$ws.Range("K105").Value = 2200
$ws.Range("M105").Value = -453

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 73.833336
$ws.Range("I7").Value = 73.833336
$ws.Range("K7").Value = 73.833336
$ws.Range("M7").Value = 39.166664
$ws.Range("H31").Value = 2456.4644
$ws.Range("J31").Value = 4221.857
$ws.Range("L31").Value = 4221.857
$ws.Range("N31").Value = -4811.857
$ws.Range("H34").Value = 2456.4644
$ws.Range("J34").Value = 4221.857
$ws.Range("L34").Value = 4221.857
$ws.Range("N34").Value = -4625.857
$ws.Range("H58").Value = 1893.5333
$ws.Range("I58").Value = 1893.5333
$ws.Range("K58").Value = 1893.5333
$ws.Range("M58").Value = -1690.5333
$ws.Range("H132").Value = 977.61536
$ws.Range("I132").Value = 977.61536
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2932.84608
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = ""
$ws.Range("H134").Value = 4135.067
$ws.Range("I134").Value = 4318
$ws.Range("K134").Value = 12954
$ws.Range("M134").Value = -10419
$ws.Range("H136").Value = 1893.5333
$ws.Range("I136").Value = 1893.5333
$ws.Range("K136").Value = 5680.5999
$ws.Range("M136").Value = -3130.5999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 999.6667
$ws.Range("J5").Value = 999.5
$ws.Range("L5").Value = 2998.5
$ws.Range("N5").Value = -3222.5
$ws.Range("H123").Value = 2165.6667
$ws.Range("I123").Value = 1748.5
$ws.Range("K123").Value = 5245.5
$ws.Range("M123").Value = -2795.5
$ws.Range("H135").Value = 999.6667
$ws.Range("J135").Value = 999.5
$ws.Range("L135").Value = 8995.5
$ws.Range("N135").Value = -14065.5
$ws.Range("H141").Value = 1386.3334
$ws.Range("I141").Value = 1386.3334
$ws.Range("K141").Value = 4159.0002
$ws.Range("M141").Value = 1020.9998

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 531.6429000000001
$ws.Range("I97").Value = 304
$ws.Range("J97").Value = 1366.3334
$ws.Range("K97").Value = 304
$ws.Range("L97").Value = 1366.3334
$ws.Range("M97").Value = 192
$ws.Range("N97").Value = -2358.3334
$ws.Range("H132").Value = 3533.25
$ws.Range("I132").Value = 3533.25
$ws.Range("K132").Value = 10599.75
$ws.Range("M132").Value = -8069.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 4
$ws.Range("I14").Value = 4
$ws.Range("K14").Value = 4
$ws.Range("M14").Value = 168
$ws.Range("H16").Value = 368.5
$ws.Range("I16").Value = 368.5
$ws.Range("K16").Value = 368.5
$ws.Range("M16").Value = -198.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 48333.332
$ws.Range("J101").Value = 48333.332
$ws.Range("L101").Value = 48333.332
$ws.Range("N101").Value = -54823.332
$ws.Range("H113").Value = 5003.9565
$ws.Range("J113").Value = 933.8570999999999
$ws.Range("L113").Value = 2801.5713
$ws.Range("N113").Value = -7141.5713
$ws.Range("H136").Value = 1902.5
$ws.Range("I136").Value = 1332.4445
$ws.Range("K136").Value = 3997.3335
$ws.Range("M136").Value = -1447.3335
